{"js": "const body = context.document.body;\n\n// Fix the typo: \"if only bulbs\" -> \"(f only bulbs\"\nconst hits = body.search(\"if only bulbs\", { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\nif (hits.items.length > 0) {\n  hits.items[0].insertText(\"(f only bulbs\", Word.InsertLocation.replace);\n}\n\n// Load every paragraph together with its inline pictures so we can\n// find the three image-only paragraphs.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  paragraph.inlinePictures.load(\"items\");\n}\nawait context.sync();\n\n// Delete each paragraph that contains an inline picture (this removes\n// both the picture and its paragraph mark, merging the remaining text\n// into a single paragraph).\nconst paragraphsToDelete = paragraphs.items.filter(\n  (paragraph) => paragraph.inlinePictures.items.length > 0\n);\nfor (const paragraph of paragraphsToDelete) {\n  paragraph.delete();\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Fix the typo first: \"if only bulbs\" -> \"(f only bulbs\"\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"if only bulbs\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"(f only bulbs\"\n$find.Execute($null, $false, $false, $false, $null, $null, $true, $null, $null, $null, 2)\n\n# Remove the three picture paragraphs entirely (the picture + its own\n# paragraph mark), leaving only the paragraph with the question text.\n# Deleting from the end keeps earlier indices stable while we work.\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.InlineShapes.Count -gt 0) {\n        $p.Range.Delete()\n    }\n}\n"}
